$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" headers on the two existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the existing sheets ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page margins used on the other sheets (0.75in/0.75in/1in/1in/0.5in/0.5in).
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Reuse the header style (bold + border + centered/top aligned) from the
# Weekly Quantity sheet's header row so the new sheet matches formatting.
$wsWeekly.Range("B1").Copy($wsForecast.Range("A1:D1"))
# Reuse the date-formatted cell style for column A (the date column).
$wsWeekly.Range("A2").Copy($wsForecast.Range("A2:A17"))

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows ---
$dates = @(44969.99999999999, 44990.99999999999, 44997.99999999999, 45032.99999999999, 45039.99999999999, 45046.99999999999, 45067.99999999999, 45081.99999999999, 45088.99999999999, 45095.99999999999, 45102.99999999999, 45109.99999999999, 45116.99999999999, 45123.99999999999, 45130.99999999999, 45137.99999999999)
$forecast = @(30, 33, 33, 37, 37, 38, 40, 42, 42, 43, 44, 44, 45, 46, 46, 47)
$yhatLower = @(-7.772459823169845, -5.645803548387163, -5.600943435029292, -3.42618063692758, -2.763727274442256, -0.8830513601029132, 0.4615694271406963, -0.2968642750519863, 4.872583302021973, 4.548188567941336, 3.933281766994537, 6.686754638595348, 8.931728969343855, 8.938732684477547, 10.70127671662891, 7.929487025140252)
$yhatUpper = @(70.68966677215519, 70.24760530526603, 72.14505758610544, 75.23559667160839, 76.63250923335282, 78.0980241616833, 79.08048056182975, 78.78778844252633, 82.65419442252319, 80.55707280048252, 82.18688686950286, 80.5050910070847, 82.45988711127045, 85.14655171706994, 86.21264715973301, 83.46172421707981)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 1).Value = $dates[$i]
    $wsForecast.Cells.Item($row, 2).Value = $forecast[$i]
    $wsForecast.Cells.Item($row, 3).Value = $yhatLower[$i]
    $wsForecast.Cells.Item($row, 4).Value = $yhatUpper[$i]
}

# Restore the original active sheet/selection (Weekly Quantity, cell A1).
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
